$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "GEN0810-Dr. Yvonne Treutel PhD-Hall 1`nGEN1809-sec-Hall 2`nCIE4818-Darwin Nicolas-Hall 5"
$ws.Range("C2").Value = "GEN2810-Mr. Howard Willms II-Hall 1`nGEN0810-Dr. Yvonne Treutel PhD-Hall 2`nGEN0809-Darwin Nicolas-Hall 5"
$ws.Range("D2").Value = "CIE1803-lab-Hall 1`nGEN1801-Prof. Eladio Franecki-Hall 2`nCIE1803-Vernie Sporer-Hall 5"
$ws.Range("E2").Value = "CIE1808-Destinee Feest-Hall 1`nCIE1808-Destinee Feest-Hall 2`nGEN0801-Darwin Nicolas-Hall 5"
$ws.Range("F2").Value = "POW1804-sec-Hall 1`nMEC0811-sec-Hall 2`nCIE2802-Destinee Feest-Hall 5"

$ws.Range("B3").Value = "GEN0806-sec-Hall 1`nGEN0807-sec-Hall 2`nGEN0806-Osvaldo Boyle PhD-Hall 5"
$ws.Range("C3").Value = "GEN1801-sec-Hall 1`nGEN1809-Laisha Schultz-Hall 2`nGEN0807-Osvaldo Boyle PhD-Hall 5"
$ws.Range("D3").Value = "CIE3804-lab-Hall 1`nGEN0802-Osvaldo Boyle PhD-Hall 2`nCIE2802-sec-Hall 5"
$ws.Range("E3").Value = "GEN1805-Mr. Howard Willms II-Hall 1`nCIE4818-sec-Hall 2`nGEN0801-sec-Hall 5"
$ws.Range("F3").Value = "CIE3804-Vernie Sporer-Hall 1`nGEN1801-Prof. Eladio Franecki-Hall 2`nGEN0801-Mr. Howard Willms II-Hall 5"

$ws.Range("B4").Value = "MEC0811-Dr. Yvonne Treutel PhD-Hall 1`nGEN0802-lab-Hall 2`nCIE1808-sec-Hall 5"
$ws.Range("C4").Value = "POW1804-Dr. Yvonne Treutel PhD-Hall 1`nPOW1804-Dr. Yvonne Treutel PhD-Hall 2"
$ws.Range("D4").Value = ""

$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(4).AutoFit()

$wb.Save()
